# Adds the 2025 squads for LAS PALMAS and REAL SOCIEDAD to Hoja1 (rows 578-635),
# appending 29 players each (Equipo / Posicion / Jugador / Precio).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mint the two new shared strings in the same order the original workbook used
# (REAL SOCIEDAD, then LAS PALMAS) via scratch cells far outside the used range,
# then clear them once the real rows below hold their own references.
$ws.Cells.Item(1000, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(1001, 1).Value = "LAS PALMAS"

$ws.Cells.Item(578, 1).Value = "LAS PALMAS"
$ws.Cells.Item(578, 2).Value = "PORTERO"
$ws.Cells.Item(578, 3).Value = "JASPER CILLESSEN"
$ws.Cells.Item(578, 4).Value = 14
$ws.Cells.Item(579, 1).Value = "LAS PALMAS"
$ws.Cells.Item(579, 2).Value = "PORTERO"
$ws.Cells.Item(579, 3).Value = "ÁLVARO VALLÉS"
$ws.Cells.Item(579, 4).Value = 13
$ws.Cells.Item(580, 1).Value = "LAS PALMAS"
$ws.Cells.Item(580, 2).Value = "PORTERO"
$ws.Cells.Item(580, 3).Value = "DINKO HORKAS"
$ws.Cells.Item(580, 4).Value = 11
$ws.Cells.Item(581, 1).Value = "LAS PALMAS"
$ws.Cells.Item(581, 2).Value = "DEFENSA"
$ws.Cells.Item(581, 3).Value = "MARVIN PARK"
$ws.Cells.Item(581, 4).Value = 11
$ws.Cells.Item(582, 1).Value = "LAS PALMAS"
$ws.Cells.Item(582, 2).Value = "DEFENSA"
$ws.Cells.Item(582, 3).Value = "SCOTT MCKENNA"
$ws.Cells.Item(582, 4).Value = 12
$ws.Cells.Item(583, 1).Value = "LAS PALMAS"
$ws.Cells.Item(583, 2).Value = "DEFENSA"
$ws.Cells.Item(583, 3).Value = "ÁLEX SUÁREZ"
$ws.Cells.Item(583, 4).Value = 12
$ws.Cells.Item(584, 1).Value = "LAS PALMAS"
$ws.Cells.Item(584, 2).Value = "DEFENSA"
$ws.Cells.Item(584, 3).Value = "DALEY SINKGRAVEN"
$ws.Cells.Item(584, 4).Value = 9
$ws.Cells.Item(585, 1).Value = "LAS PALMAS"
$ws.Cells.Item(585, 2).Value = "DEFENSA"
$ws.Cells.Item(585, 3).Value = "ÁLEX MUÑOZ"
$ws.Cells.Item(585, 4).Value = 12
$ws.Cells.Item(586, 1).Value = "LAS PALMAS"
$ws.Cells.Item(586, 2).Value = "DEFENSA"
$ws.Cells.Item(586, 3).Value = "MIKA MÁRMOL"
$ws.Cells.Item(586, 4).Value = 13
$ws.Cells.Item(587, 1).Value = "LAS PALMAS"
$ws.Cells.Item(587, 2).Value = "DEFENSA"
$ws.Cells.Item(587, 3).Value = "VALENTÍN PEZZOLESI"
$ws.Cells.Item(587, 4).Value = 9
$ws.Cells.Item(588, 1).Value = "LAS PALMAS"
$ws.Cells.Item(588, 2).Value = "DEFENSA"
$ws.Cells.Item(588, 3).Value = "JUAN HERZOG"
$ws.Cells.Item(588, 4).Value = 11
$ws.Cells.Item(589, 1).Value = "LAS PALMAS"
$ws.Cells.Item(589, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(589, 3).Value = "FABIO GONZÁLEZ"
$ws.Cells.Item(589, 4).Value = 10
$ws.Cells.Item(590, 1).Value = "LAS PALMAS"
$ws.Cells.Item(590, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(590, 3).Value = "JAVI MUÑOZ"
$ws.Cells.Item(590, 4).Value = 12
$ws.Cells.Item(591, 1).Value = "LAS PALMAS"
$ws.Cells.Item(591, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(591, 3).Value = "PEJIÑO"
$ws.Cells.Item(591, 4).Value = 11
$ws.Cells.Item(592, 1).Value = "LAS PALMAS"
$ws.Cells.Item(592, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(592, 3).Value = "ALBERTO MOLEIRO"
$ws.Cells.Item(592, 4).Value = 14
$ws.Cells.Item(593, 1).Value = "LAS PALMAS"
$ws.Cells.Item(593, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(593, 3).Value = "BENITO RAMÍREZ"
$ws.Cells.Item(593, 4).Value = 12
$ws.Cells.Item(594, 1).Value = "LAS PALMAS"
$ws.Cells.Item(594, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(594, 3).Value = "ENZO LOIODICE"
$ws.Cells.Item(594, 4).Value = 12
$ws.Cells.Item(595, 1).Value = "LAS PALMAS"
$ws.Cells.Item(595, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(595, 3).Value = "KIRIAN RODRÍGUEZ"
$ws.Cells.Item(595, 4).Value = 14
$ws.Cells.Item(596, 1).Value = "LAS PALMAS"
$ws.Cells.Item(596, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(596, 3).Value = "JOSÉ CAMPAÑA"
$ws.Cells.Item(596, 4).Value = 12
$ws.Cells.Item(597, 1).Value = "LAS PALMAS"
$ws.Cells.Item(597, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(597, 3).Value = "MANU FUSTER"
$ws.Cells.Item(597, 4).Value = 13
$ws.Cells.Item(598, 1).Value = "LAS PALMAS"
$ws.Cells.Item(598, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(598, 3).Value = "VITI ROZADA"
$ws.Cells.Item(598, 4).Value = 11
$ws.Cells.Item(599, 1).Value = "LAS PALMAS"
$ws.Cells.Item(599, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(599, 3).Value = "IVÁN GIL"
$ws.Cells.Item(599, 4).Value = 10
$ws.Cells.Item(600, 1).Value = "LAS PALMAS"
$ws.Cells.Item(600, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(600, 3).Value = "ADNAN JANUZAJ"
$ws.Cells.Item(600, 4).Value = 13
$ws.Cells.Item(601, 1).Value = "LAS PALMAS"
$ws.Cells.Item(601, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(601, 3).Value = "DÁRIO ESSUGO"
$ws.Cells.Item(601, 4).Value = 12
$ws.Cells.Item(602, 1).Value = "LAS PALMAS"
$ws.Cells.Item(602, 2).Value = "DELANTERO"
$ws.Cells.Item(602, 3).Value = "OLIVER MCBURNIE"
$ws.Cells.Item(602, 4).Value = 13
$ws.Cells.Item(603, 1).Value = "LAS PALMAS"
$ws.Cells.Item(603, 2).Value = "DELANTERO"
$ws.Cells.Item(603, 3).Value = "MARC CARDONA"
$ws.Cells.Item(603, 4).Value = 12
$ws.Cells.Item(604, 1).Value = "LAS PALMAS"
$ws.Cells.Item(604, 2).Value = "DELANTERO"
$ws.Cells.Item(604, 3).Value = "JAIME MATA"
$ws.Cells.Item(604, 4).Value = 12
$ws.Cells.Item(605, 1).Value = "LAS PALMAS"
$ws.Cells.Item(605, 2).Value = "DELANTERO"
$ws.Cells.Item(605, 3).Value = "SANDRO RAMÍREZ"
$ws.Cells.Item(605, 4).Value = 13
$ws.Cells.Item(606, 1).Value = "LAS PALMAS"
$ws.Cells.Item(606, 2).Value = "DELANTERO"
$ws.Cells.Item(606, 3).Value = "FÁBIO SILVA"
$ws.Cells.Item(606, 4).Value = 12
$ws.Cells.Item(607, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(607, 2).Value = "PORTERO"
$ws.Cells.Item(607, 3).Value = "ÁLEX REMIRO"
$ws.Cells.Item(607, 4).Value = 17
$ws.Cells.Item(608, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(608, 2).Value = "PORTERO"
$ws.Cells.Item(608, 3).Value = "AITOR FRAGA"
$ws.Cells.Item(608, 4).Value = 14
$ws.Cells.Item(609, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(609, 2).Value = "PORTERO"
$ws.Cells.Item(609, 3).Value = "UNAI MARRERO"
$ws.Cells.Item(609, 4).Value = 15
$ws.Cells.Item(610, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(610, 2).Value = "DEFENSA"
$ws.Cells.Item(610, 3).Value = "ÁLVARO ODRIOZOLA"
$ws.Cells.Item(610, 4).Value = 13
$ws.Cells.Item(611, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(611, 2).Value = "DEFENSA"
$ws.Cells.Item(611, 3).Value = "AIHEN MUÑOZ"
$ws.Cells.Item(611, 4).Value = 13
$ws.Cells.Item(612, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(612, 2).Value = "DEFENSA"
$ws.Cells.Item(612, 3).Value = "ARITZ ELUSTONDO"
$ws.Cells.Item(612, 4).Value = 15
$ws.Cells.Item(613, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(613, 2).Value = "DEFENSA"
$ws.Cells.Item(613, 3).Value = "JAVI LÓPEZ"
$ws.Cells.Item(613, 4).Value = 14
$ws.Cells.Item(614, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(614, 2).Value = "DEFENSA"
$ws.Cells.Item(614, 3).Value = "HAMARI TRAORÉ"
$ws.Cells.Item(614, 4).Value = 14
$ws.Cells.Item(615, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(615, 2).Value = "DEFENSA"
$ws.Cells.Item(615, 3).Value = "JON PACHECO"
$ws.Cells.Item(615, 4).Value = 15
$ws.Cells.Item(616, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(616, 2).Value = "DEFENSA"
$ws.Cells.Item(616, 3).Value = "IGOR ZUBELDIA"
$ws.Cells.Item(616, 4).Value = 17
$ws.Cells.Item(617, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(617, 2).Value = "DEFENSA"
$ws.Cells.Item(617, 3).Value = "SERGIO GÓMEZ"
$ws.Cells.Item(617, 4).Value = 17
$ws.Cells.Item(618, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(618, 2).Value = "DEFENSA"
$ws.Cells.Item(618, 3).Value = "NAYEF AGUERD"
$ws.Cells.Item(618, 4).Value = 16
$ws.Cells.Item(619, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(619, 2).Value = "DEFENSA"
$ws.Cells.Item(619, 3).Value = "JON ARAMBURU"
$ws.Cells.Item(619, 4).Value = 13
$ws.Cells.Item(620, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(620, 2).Value = "DEFENSA"
$ws.Cells.Item(620, 3).Value = "JON MARTÍN"
$ws.Cells.Item(620, 4).Value = 12
$ws.Cells.Item(621, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(621, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(621, 3).Value = "MARTÍN ZUBIMENDI"
$ws.Cells.Item(621, 4).Value = 18
$ws.Cells.Item(622, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(622, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(622, 3).Value = "URKO GONZÁLEZ"
$ws.Cells.Item(622, 4).Value = 13
$ws.Cells.Item(623, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(623, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(623, 3).Value = "JON OLASAGASTI"
$ws.Cells.Item(623, 4).Value = 15
$ws.Cells.Item(624, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(624, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(624, 3).Value = "ARSEN ZAKHARYAN"
$ws.Cells.Item(624, 4).Value = 16
$ws.Cells.Item(625, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(625, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(625, 3).Value = "JON MAGUNAZELAIA"
$ws.Cells.Item(625, 4).Value = 14
$ws.Cells.Item(626, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(626, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(626, 3).Value = "LUKA SUCIC"
$ws.Cells.Item(626, 4).Value = 18
$ws.Cells.Item(627, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(627, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(627, 3).Value = "BEÑAT TURRIENTES"
$ws.Cells.Item(627, 4).Value = 17
$ws.Cells.Item(628, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(628, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(628, 3).Value = "BRAIS MÉNDEZ"
$ws.Cells.Item(628, 4).Value = 18
$ws.Cells.Item(629, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(629, 2).Value = "MEDIOCENTRO"
$ws.Cells.Item(629, 3).Value = "PABLO MARÍN"
$ws.Cells.Item(629, 4).Value = 14
$ws.Cells.Item(630, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(630, 2).Value = "DELANTERO"
$ws.Cells.Item(630, 3).Value = "SHERALDO BECKER"
$ws.Cells.Item(630, 4).Value = 15
$ws.Cells.Item(631, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(631, 2).Value = "DELANTERO"
$ws.Cells.Item(631, 3).Value = "ANDER BARRENETXEA"
$ws.Cells.Item(631, 4).Value = 15
$ws.Cells.Item(632, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(632, 2).Value = "DELANTERO"
$ws.Cells.Item(632, 3).Value = "ORRI ÓSKARSSON"
$ws.Cells.Item(632, 4).Value = 14
$ws.Cells.Item(633, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(633, 2).Value = "DELANTERO"
$ws.Cells.Item(633, 3).Value = "MIKEL OYARZABAL"
$ws.Cells.Item(633, 4).Value = 19
$ws.Cells.Item(634, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(634, 2).Value = "DELANTERO"
$ws.Cells.Item(634, 3).Value = "UMAR SADIQ"
$ws.Cells.Item(634, 4).Value = 17
$ws.Cells.Item(635, 1).Value = "REAL SOCIEDAD"
$ws.Cells.Item(635, 2).Value = "DELANTERO"
$ws.Cells.Item(635, 3).Value = "TAKEFUSA KUBO"
$ws.Cells.Item(635, 4).Value = 20

# Drop the scratch cells now that they are no longer the sole reference
$ws.Cells.Item(1000, 1).Value = ""
$ws.Cells.Item(1001, 1).Value = ""

# Match the saved selection / active cell from the source workbook
$ws.Activate()
$ws.Range("A607:A635").Select()

Write-Output "Added LAS PALMAS and REAL SOCIEDAD rows (578-635)"
